$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.052.92'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.494.65'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.64'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.72'
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.517.76'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  -2.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.33'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.347'
$ws.Range('E13').Value = '  -3.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.968.87'
$ws.Range('E14').Value = '  +0.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.00'
$ws.Range('E15').Value = '  -2.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '59.020.99'
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.521.57'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.11'
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.28'
$ws.Range('E20').Value = '  -0.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.50'
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('E23').Value = '  +1.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.47'
$ws.Range('E24').Value = '  +2.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.420'
$ws.Range('E25').Value = '  -1.86%  '
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.59'
$ws.Range('E28').Value = '  -2.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.71'
$ws.Range('E29').Value = '  -3.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0766'
$ws.Range('E30').Value = '  -1.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.77'
$ws.Range('E31').Value = '  -0.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '165.84'
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.17'
$ws.Range('E33').Value = '  +5.29%  '
$ws.Range('E34').Value = '  +1.92%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.42'
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.10'
$ws.Range('E37').Value = '  -3.92%  '
$ws.Range('E38').Value = '  -3.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.62'
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.812'
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.62'
$ws.Range('E41').Value = '  -2.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '285.06'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.16'
$ws.Range('E43').Value = '  -1.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '132.66'
$ws.Range('E44').Value = '  +7.97%  '
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.89'
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0928'
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0507'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0220'
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.24'
$ws.Range('E51').Value = '  -3.20%  '
